$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, shifting existing rows 10-14 down to 11-15
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with data for a newer weekly report entry
$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(10, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(10, 4).Value = 44673
$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 5).Value = 15
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100101
$ws.Cells.Item(10, 8).Value = "Berries"
$ws.Cells.Item(10, 9).Value = 100101007
$ws.Cells.Item(10, 10).Value = "Kiwi"
$ws.Cells.Item(10, 11).Value = "Hayward"
$ws.Cells.Item(10, 12).Value = "Especial"
$ws.Cells.Item(10, 13).Value = 400
$ws.Cells.Item(10, 14).Value = 14000
$ws.Cells.Item(10, 15).Value = 15000
$ws.Cells.Item(10, 16).Value = 14500
$ws.Cells.Item(10, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(10, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(10, 19).Value = 1450
$ws.Cells.Item(10, 20).Value = 10
